$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A helper/scratch cell used to force text (string) values for the "weight"
# column without leaving a residual text NumberFormat style on the real
# cells (Excel auto-converts numeric-looking strings like "0.25" to a
# number unless pasted in from a cell that is already formatted as text).
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"

# Step 1: update the existing 40 answer rows (rows 2-41):
#   - column B gets an extra leading space (3 spaces -> 4 spaces), making
#     room for the student name captured alongside the answer key;
#   - column C switches from the literal label "Peso" to the actual
#     weight value of each question (0.25), kept as text.
for ($r = 2; $r -le 41; $r++) {
    $letter = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 2).Value2 = " " + $letter

    $helper.Value2 = "0.25"
    $helper.Copy()
    $ws.Cells.Item($r, 3).PasteSpecial(-4163) | Out-Null
}
$helper.Clear() | Out-Null
$excel.CutCopyMode = 0

# Step 2: duplicate the 40-question block (rows 2-41) into rows 42-81,
# preserving formatting/styles, to extend the gabarito (answer key) with a
# second set of questions.
$src = $ws.Range("A2:C41")
$dst = $ws.Range("A42")
$src.Copy($dst)
$excel.CutCopyMode = 0

# Step 3: renumber the question index in column A for the appended block
# so numbering continues (40-79) instead of repeating (0-39).
for ($r = 42; $r -le 81; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $r - 2
}
